$wb = $excel.ActiveWorkbook
$permisos = $wb.Worksheets.Item("Permisos")

# Populate row 15 on Permisos sheet (model / permission-type / group) which
# was previously left blank, causing the Datos!10 formulas to render
# placeholder values like "access__for_group_name_,...".
$permisos.Range("B15").Value = "certification_contract"
$permisos.Range("C15").Value = "read_write_create_unlink"
$permisos.Range("D15").Value = "certifications_administrator"

# Recalculate so cached formula results (Permisos!E15, Datos!A10:H10, etc.)
# reflect the new inputs.
$excel.Calculate()

# Update the active selection to match the cell that was edited.
$permisos.Activate()
$permisos.Range("E15").Select()
